# Updated diary for thuc nguyen
# Rows 19 and 20 were placeholder template rows ("<what day?>", "<what time?>", ...);
# they are turned into real diary entries, matching the look/format of row 18
# (date value + normal-entry styling) while row 21 is left as the remaining
# placeholder template row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the already-filled-in row 18 down onto rows 19-20 so
# that they pick up the "real entry" styling (date number format on column A,
# wrapped/top-aligned text on B-G, etc.) instead of the plain placeholder style.
$ws.Range("A18:G18").Copy()
$ws.Range("A19:G20").PasteSpecial(-4122)

# Row 19
$ws.Range("A19").Value = 43874.0
$ws.Range("B19").Value = "5:00pm - 8:00pm"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = "Review the class survey`nUnderstand the next key expert practices`nWhat are stakeholders? The functionality of our system? The key developers?"
$ws.Range("E19").Value = "1) Future topics to cover in the class, feedback on the most and least useful topics. `n2) Key Expert Practices`n    #7. Prioritize among stakeholders`n    #8. Move along levels of abstraction`n    #9. Do something else`n3) Stakeholders, essential functional and essential non-functional aspects of the system, and key developers such as core maintainers, team members, developers, testers, triagers, documentation writers"
$ws.Range("F19").Value = "It was interesting to step away from code and look at other aspects of our system, in a higher level view. I felt that we learned how to read code very well from the past lectures, and that we were able to find the core or essence of our system, but not about who made this project, and the people that were interested in its growth."
$ws.Range("G19").Value = "Satisfied, and enlightened to take pm a different perspective about the system."

# Row 20
$ws.Range("A20").Value = 43880.0
$ws.Range("B20").Value = "8:00pm - 11:00pm"
$ws.Range("C20").Value = "Deon and Harry"
$ws.Range("D20").Value = "Learn and share our research on the stakeholders, functionality, key developers, and the issues of our system"
$ws.Range("E20").Value = "We found and described the stakeholders, functionality, key developers, and 5 issues that we can potentially solve. "
$ws.Range("F20").Value = "It was surprising to find that our system was solely lead by a single developer, Adam. He does have a team of core maintainers however. Our system also had more stakeholders than I had expected, including the people of Venezuela how play the game even today, in order to convert in-game currency to real money because their economy's financial infrastructure is hyperinflated. The way we went about searching for this information was very straightforward. We looked through Github, the Runelite website, Jagex's website, google, and so on."
$ws.Range("G20").Value = "Surprised especially because it was the first time I read that a game and client had an impact on people's livelihoods."

# Keep the original, compact row heights (auto-fit would otherwise grow them
# to fit the long wrapped text we just entered).
$ws.Rows.Item(19).RowHeight = 14.25
$ws.Rows.Item(20).RowHeight = 14.25

# Row 21 stays the remaining unfilled placeholder/template row, untouched.
